$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format temporarily so numeric-looking strings
# (e.g. "1.01") are stored as text, matching the source inlineStr cells.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "25.594.00"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "1.624.08"
$ws.Range("E3").Value = "  -1.21%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "213.83"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "0.498"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").Value = "1.01"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "0.0633"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.623.64"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.22"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.853.23"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "0.549"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "0.0₃0759"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "62.39"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "25.631.80"
$ws.Range("E18").Value = "  -1.12%  "
$ws.Range("D19").Value = "1.01"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "4.39"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").Value = "192.62"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "9.88"
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("D23").Value = "6.17"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "1.77"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").Value = "139.54"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "0.119"
$ws.Range("E27").Value = "  -3.75%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "15.36"
$ws.Range("E29").Value = "  -0.91%  "
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").Value = "0.0483"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("D34").Value = "1.58"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "0.890"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").Value = "2.55"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "0.541"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").Value = "1.101.79"
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("D40").Value = "0.0154"
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "5.55"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("D43").Value = "99.58"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").Value = "0.792"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "1.759.87"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("E46").Value = "  -4.79%  "
$ws.Range("D47").Value = "54.67"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.418"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "2.36"
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.0500"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.57"
$ws.Range("E51").Value = "  -2.59%  "

# Restore normal style on column D so no stray number-format styling remains
$colD.Style = "Normal"

Write-Output "done"
